# AWS_Progress_Tracker.xlsx - "added code for new services"
# Adds 4 new tracker rows (23-26: ELB, GuardDuty, Secrets Manager, Athena)
# to the bottom of the previously-blank row stubs 23-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 23 - ELB
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = 44631
$ws.Range("A23").NumberFormat = "m/d/yy"
$ws.Range("B23").Value = 44631
$ws.Range("B23").NumberFormat = "m/d/yy"
$ws.Range("C23").Value = "ELB"
$ws.Range("D23").Value = "Automatically distribute incoming application traffic in one or more availability zones`nMonitor health and performance of applications"
$ws.Range("F23").Value = "Tried to create a load balancer and attached it to t2.micro EC2 instance but failed"
$ws.Range("J23").Value = "EC2`nLambda`nFargate`nEKS`nECS`nWAF`nCertificate Manager`nCognito"

# ---------------------------------------------------------------------------
# Row 24 - GuardDuty
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = 44631
$ws.Range("A24").NumberFormat = "m/d/yy"
$ws.Range("B24").Value = 44631
$ws.Range("B24").NumberFormat = "m/d/yy"
$ws.Range("C24").Value = "GuardDuty"
$ws.Range("D24").Value = "Threat Detection`nMonitor malicious activity"
$ws.Range("F24").Value = "Created GuardDuty Detector resource using terraform"
$ws.Range("H24").Value = "CloudWatch`nCloudTrail`nS3"
$ws.Range("J24").Value = "S3"

# ---------------------------------------------------------------------------
# Row 25 - Secrets Manager
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = 44631
$ws.Range("A25").NumberFormat = "m/d/yy"
$ws.Range("B25").Value = 44631
$ws.Range("B25").NumberFormat = "m/d/yy"
$ws.Range("C25").Value = "Secrets Manager"
$ws.Range("D25").Value = "Rotate, manage, and retrieve database credentials, API keys`nSecret encryption at rest"
$ws.Range("E25").Value = "Not available for free tier access"
$ws.Range("F25").Value = "Not available for free tier access"
$ws.Range("G25").Value = "Not available for free tier access"
$ws.Range("H25").Value = "CloudTrail"
$ws.Range("J25").Value = "RDS`nRedshift `nDocumentDB`nKMS`n"

# ---------------------------------------------------------------------------
# Row 26 - Athena
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = 44631
$ws.Range("A26").NumberFormat = "m/d/yy"
$ws.Range("B26").Value = 44631
$ws.Range("B26").NumberFormat = "m/d/yy"
$ws.Range("C26").Value = "Athena"
$ws.Range("D26").Value = "Interactive query service to analyze data in S3 using standard SQL`nServerless"
$ws.Range("E26").Value = "Not available for free tier access"
$ws.Range("F26").Value = "Not available for free tier access"
$ws.Range("G26").Value = "Not available for free tier access"
$ws.Range("H26").Value = "CloudTrail"
$ws.Range("J26").Value = "Glue`nS3`nQuickSight"

# ---------------------------------------------------------------------------
# Hyperlinks - added in the same order the author clicked them in, so the
# generated relationship ids line up (I23, G23, G24, I24, I25, I26).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("I23"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/AWS%20Services%20Actions.xlsx")
$ws.Hyperlinks.Add($ws.Range("G23"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/loadbalancer.tf")
$ws.Range("G23").WrapText = $true
$ws.Range("G23").VerticalAlignment = -4160

$ws.Hyperlinks.Add($ws.Range("G24"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/guardduty.tf")
$ws.Range("G24").WrapText = $true
$ws.Range("G24").VerticalAlignment = -4160

$ws.Hyperlinks.Add($ws.Range("I24"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/AWS%20Services%20Actions.xlsx")
$ws.Hyperlinks.Add($ws.Range("I25"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/AWS%20Services%20Actions.xlsx")
$ws.Hyperlinks.Add($ws.Range("I26"), "https://github.com/ShivamGautam98/AWS-terraform/blob/main/AWS%20Services%20Actions.xlsx")

# ---------------------------------------------------------------------------
# Row heights (Excel auto-fit the wrapped text to whole multiples of the
# default 14.5pt line height: 8 lines / 5 lines / 5 lines / 5 lines).
# ---------------------------------------------------------------------------
$ws.Rows.Item(23).RowHeight = 116
$ws.Rows.Item(24).RowHeight = 72.5
$ws.Rows.Item(25).RowHeight = 72.5
$ws.Rows.Item(26).RowHeight = 72.5

# ---------------------------------------------------------------------------
# Selection moved to H26 as the last-touched cell.
# ---------------------------------------------------------------------------
$ws.Range("H26").Select()
